$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the contents of row 4 (A4:F4) while keeping cell styles/formatting
$ws.Range("A4:F4").ClearContents()

# Update the selection to A4 with the range A4:F4 selected
$ws.Range("A4:F4").Select()
